$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the existing "iter_cPCA" run to note that it used 10 batches.
# ---------------------------------------------------------------------------
$ws.Range("P42").Value2 = "iter_cPCA 10 batches"

# ---------------------------------------------------------------------------
# 2) The sheet already reserves rows 49-56 (styled but empty) below the
#    existing results table. Populate rows 49-52 with four new test runs,
#    copying the formatting from the row above (row 48) so the new rows pick
#    up the same cell styles used throughout the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A48:P48").Copy() | Out-Null
$ws.Range("A49:P52").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Also re-apply that same formatting (without values) to the remaining blank
# rows 53-56 so every row in the block shares one consistent set of styles.
$ws.Range("I48:P48").Copy() | Out-Null
$ws.Range("I53:P56").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Rows 49 & 50 were using the shorter "blank row" height (18.75); bring them
# up to the standard data-row height (19.5) used elsewhere in the table.
$ws.Rows.Item(49).RowHeight = 19.5
$ws.Rows.Item(50).RowHeight = 19.5

function Set-Row($r, $target, $overlap10, $overlap20, $notes) {
    $ws.Range("A$r").Value2 = "ukb51139_subset.csv"
    $ws.Range("B$r").Value2 = "28012 x 1081"
    $ws.Range("C$r").Value2 = "all"
    $ws.Range("D$r").Value2 = "no events"
    $ws.Range("E$r").Value2 = $target
    $ws.Range("F$r").Value2 = "zscore"
    $ws.Range("G$r").Value2 = "median"
    $ws.Range("H$r").Value2 = "none"
    $ws.Range("I$r").Value2 = 25
    $ws.Range("K$r").Value2 = "N/A"
    $ws.Range("L$r").Value2 = $overlap10
    $ws.Range("M$r").Value2 = $overlap20
    $ws.Range("N$r").Value2 = "N/A"
    $ws.Range("O$r").Value2 = "N/A"
    $ws.Range("P$r").Value2 = $notes
}

Set-Row 49 "> 160/100" "-18.5 & -5"   "-27.6 & -18.4"  "10 batches"
Set-Row 50 "> 160/100" "21.0 & 7.2"   "-118.8 & -51.8" "15 batches"
Set-Row 51 "> 140/80"  "-136 & -49.3" "19.6 & 16.9"    "15 batches (shuffled data)"
Set-Row 52 "> 140/80"  "-45.6 & -22.2" "4.9 & 4.2"     "25 batches"
